# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund breakdown, same layout as the
#    "2021-Q4" sheet) positioned right before the "总计" (totals) sheet.
# 2. Insert a new first data-row into the "总计" sheet summarising the
#    2022-Q1 quarter (holding count + market value), shifting the existing
#    rows down and renumbering the running index in column A.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")
$totalBefore = $wb.Worksheets.Item("总计")

# --- 1. Build the "2022-Q1" sheet from a copy of "2021-Q4" (same header
#        layout/styles), trim it down to a single data row, then fill in
#        the new figures. --------------------------------------------------
$q4.Copy($totalBefore)

# NOTE: inserting a sheet shifts everybody's position in the Worksheets
# collection, so any worksheet handle obtained *before* the Copy() call
# (e.g. $totalBefore above) is stale afterwards - always re-resolve sheets
# by name after a sheet-collection change like this.
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q1"

# Drop the four extra fund rows that came along with the copy, keeping
# only the header (row 1) and one data row (row 2).
$new.Range("A3:A6").EntireRow.Delete()

# The text-ish numeric-looking columns must stay text cells (matching the
# rest of the workbook), so force text format before writing them.
$new.Range("B2:G2").NumberFormat = "@"

$new.Range("A2").Value = 0
$new.Range("B2").Value = "008099"
$new.Range("C2").Value = "广发价值领先混合"
$new.Range("D2").Value = "61.82"
$new.Range("E2").Value = "83.88"
$new.Range("F2").Value = "4.33"
$new.Range("G2").Value = "2.6768"
$new.Range("H2").Value = 6

# Drop the format override we used to coerce text cells, so the data row
# ends up unstyled like the other sheets' data rows.
$new.Range("B2:G2").ClearFormats()

# --- 2. Insert the 2022-Q1 summary row at the top of "总计". ---------------
# Re-resolve "总计" by name since the Copy() above moved it one slot over.
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# The freshly inserted row has no format on column A; pull the running
# index style from the row below (which used to be row 2) before we
# touch the values.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 2.68

# Renumber the running index (column A) for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
